$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in the i18n placeholder used as the "tickets in work" column header:
# "{d.i18n.tocketsInWork}" -> "{d.i18n.ticketsInWork}"
$ws.Range("E1").Value = "{d.i18n.ticketsInWork}"
